# MS XL Addin updates in function names + documentation
#
# Renames custom Excel add-in functions referenced in labels/formulas:
#   SLTAX.DETAILS               -> SLTAX.CALCULATE_DETAILS
#   SLTAX.WITHPREMIUM           -> SLTAX.CALCULATE_WITHPREMIUM
#   SLTAX.RATESDETAILS          -> SLTAX.RATES_DETAILS
#   SLTAX.HISTORICALRATEDETAILS -> SLTAX.HISTORICALRATE_DETAILS

$wb = $excel.ActiveWorkbook

# ----- Calculator sheet -----
$ws = $wb.Worksheets.Item("Calculator")

$ws.Range("A11").Value = "Detailed Breakdown (SLTAX.CALCULATE_DETAILS)"
$ws.Range("C13").Formula = "=SLTAX.CALCULATE_DETAILS(A13, B13)"
$ws.Range("C14").Formula = "=SLTAX.CALCULATE_DETAILS(A14, B14)"

$ws.Range("A17").Value = "Compact View (SLTAX.CALCULATE_WITHPREMIUM)"
$ws.Range("B19").Formula = "=SLTAX.CALCULATE_WITHPREMIUM(""Florida"", 15000)"

# ----- Detailed Rates sheet -----
$ws = $wb.Worksheets.Item("Detailed Rates")

$ws.Range("A1").Value = "Complete Rate Details (SLTAX.RATES_DETAILS)"
$ws.Range("A4").Formula = "=SLTAX.RATES_DETAILS()"

# ----- Historical sheet -----
$ws = $wb.Worksheets.Item("Historical")

$ws.Range("A11").Formula = "=SLTAX.HISTORICALRATE_DETAILS(""Texas"", ""2024-01-01"")"
$ws.Range("A15").Formula = "=SLTAX.HISTORICALRATE_DETAILS(""Texas"", ""2024-01-01"", TRUE)"

# ----- Quick Reference sheet -----
$ws = $wb.Worksheets.Item("Quick Reference")

$ws.Range("A5").Value = "SLTAX.CALCULATE_DETAILS(state, premium, [multiline])"
$ws.Range("C5").Value = "SLTAX.CALCULATE_DETAILS(""CA"", 10000)"

$ws.Range("A6").Value = "SLTAX.CALCULATE_WITHPREMIUM(state, premium)"
$ws.Range("C6").Value = "SLTAX.CALCULATE_WITHPREMIUM(""FL"", 15000)"

$ws.Range("A10").Value = "SLTAX.RATES_DETAILS()"
$ws.Range("C10").Value = "SLTAX.RATES_DETAILS()"

$ws.Range("A12").Value = "SLTAX.HISTORICALRATE_DETAILS(state, date, [multiline])"
$ws.Range("C12").Value = "SLTAX.HISTORICALRATE_DETAILS(""TX"", ""2024-01-01"")"
